$d = $word.ActiveDocument

function Get-ParaByText($doc, $needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

function Add-ItalicParaAfter($doc, $para, $text) {
    $para.Range.InsertParagraphAfter()
    $newPara = $para.Next()
    $newRange = $newPara.Range
    $newRange.Text = $text
    $italicRange = $doc.Range($newRange.Start, $newRange.Start + $text.Length)
    $italicRange.Font.Italic = $true
}

# 1. Update activation date
$ok1 = $d.Content.Find.Execute("Ativação: 01/01/2020", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2023", 2)
Write-Output "step1 (activation date): $ok1"

# 2. Add English translation of "Objetivos" paragraph (italic)
$pObjetivos = Get-ParaByText $d "Apresentar o formalismo"
Add-ItalicParaAfter $d $pObjetivos "To present the formalism for the description of quantum systems. Study several applications of the time-independent Schroedinger equation. Describe the electronic structure of atoms and molecules."
Write-Output "step2 (objectives EN paragraph): done"

# 3. Remove two instructors from the "Docente(s) Responsável(eis)" list
$ok3 = $d.Content.Find.Execute("6279110 - Carlos Alberto Moreira dos Santos^l6495737 - Durval Rodrigues Junior^l", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
Write-Output "step3 (remove two instructors): $ok3"

# 4. Add English translation of "Programa resumido" paragraph (italic)
$pResumido = Get-ParaByText $d "Introdução aos conceitos da Mecânica Quântica"
Add-ItalicParaAfter $d $pResumido "• Introduction to the concepts of Quantum Mechanics. • Mathematical tools of Quantum Mechanics. • The Schrödinger equation and one- and three-dimensional applications. • Quantum formalism. • Problems in rectangular coordinates and spherical coordinates. • Hydrogen atoms and orbitals. • General properties of angular momentum. • Spin. • Fermions and bosons."
Write-Output "step4 (summary EN paragraph): done"

# 5. Replace the "Programa" paragraph text with the new Portuguese syllabus
$oldPrograma = "Origens das ideias fundamentais da Mecânica Quântica.Dualidade onda partícula. Principio de Heisenberg.• Os postulados da Mecânica Quântica. Ferramentas matemáticas da Mecânica Quântica. O Espaço de Hilbert e a Equação de Onda. Notação de Dirac. Operadores e Bases. Representação matricial. A equação de Schroedinger e aplicações unidimensionais. Barreira de potencial. Poço de potencial. Oscilador harmônico. Problemas tridimensionais. Problemas em coordenadas retangulares. Problemas em coordenadas esféricas. Átomos com um elétron. Teoria geral.• Momento angular. Spin do elétron.• Propriedades gerais do momento angular"
$newPrograma = "• Equação de Schrödinger. • Função de onda e interpretação estatística da mecânica quântica. • Valores esperados e operadores. Os operadores posição e momento; operadores energia cinética e potencial; o operador Hamiltoniano. • A equação de Schrödinger independente do tempo. Separação de variáveis e estados estacionários. • Aplicações unidimensionais:  poço quadrado infinito; oscilador harmônico; partícula livre;  transformada de Fourier e sua relação com o princípio da incerteza de Heisenberg; Poços e barreiras de potencial. • Formalismo quântico: opserváveis e operadores hermitianos. Estados determinados, autoestados e autovalores de operadores hermitianos. Base de autoestados; interpretação estatística generalizada: medidas de observáveis e suas probabilidades. Comutadores e operadores que compartilham autoestados ; princípio da incerteza generalizado. • Mecânica Quântica em três dimensões. • Átomo de hidrogênio: modelo de Bohr e o número quântico principal. Solução completa e os demais números quânticos. • Coordenadas esféricas e Momento angular.  • Momento angulas de spin. • Problemas de muitos corpos. • Partículas idênticas: férmions e bósons."
$ok5 = $d.Content.Find.Execute($oldPrograma, $true, $false, $false, $false, $false, $true, 1, $false, $newPrograma, 2)
Write-Output "step5 (programa PT rewrite): $ok5"

# 6. Add English translation of the (now updated) "Programa" paragraph (italic)
$pPrograma = Get-ParaByText $d "Equação de Schrödinger. • Função de onda"
Add-ItalicParaAfter $d $pPrograma "• Schrödinger's equation. • Wave function and statistical interpretation of quantum mechanics. • Expected values and operators. The position and moment operators; kinetic and potential energy operators; the Hamiltonian operator. • The time-independent Schrödinger equation. Separation of variables and steady states. • One-dimensional applications: infinite square well; harmonic oscillator; free particle; Fourier transform and its relationship with the Heisenberg uncertainty principle; Potential square wells and barriers. • Quantum formalism: hermitian operators and observables. Determined states, eigenstates and eigenvalues of Hermitian operators. Basis of Eigenstates; generalized statistical interpretation: measures of observables and their probabilities. Comutators and operators that share eigenstates; generalized uncertainty principle. • Quantum Mechanics in three dimensions. • Hydrogen atom: Bohr model and the principal quantum number. Complete solution and the other quantum numbers. • Spherical coordinates and Angular momentum. • Spin angular momentum. • Many-body problems. • Identical particles: fermions and bosons."
Write-Output "step6 (programa EN paragraph): done"

# 7. Update the evaluation criterion text
$ok7 = $d.Content.Find.Execute("Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2.", $true, $false, $false, $false, $false, $true, 1, $false, "Média aritmética de três provas: P1 (peso 1), P2 (peso 1) e P3 (peso 2).", 2)
Write-Output "step7 (grading criterion): $ok7"

Write-Output "done"
